$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F12").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("F28").Value = 0
$ws.Range("F31").Value = -3
$ws.Range("F36").Value = -1
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 3
$ws.Range("F46").Value = -3
$ws.Range("F49").Value = 1
$ws.Range("F50").Value = -1
$ws.Range("F51").Value = 2
$ws.Range("F54").Value = 1
$ws.Range("F59").Value = 2
$ws.Range("F68").Value = 4
$ws.Range("F71").Value = 5
